$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$v = $ws.Range("B2").Value()
$v2 = $v.Replace("Wed, 01 Jan 2020", "Thu, 02 Jan 2020")
$v3 = $v2.Replace("23:19:04", "20:48:57")
$ws.Range("B2").Value = $v3
$check = $ws.Range("B2").Value()
Write-Host $check.Substring(0,400)
